# Updated cryptos list on Sat Sep 21 05:30:30 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the Stacks / FirstDigitalUSD rows (40/41) back into rank order.
#
# Note: several Price values are plain decimals (e.g. "573.33") that Excel's
# Range.Value setter would otherwise auto-coerce to floating-point numbers,
# losing the exact textual formatting the source data uses. A leading
# apostrophe is prefixed for those (Excel's standard "treat as text"
# quote-prefix convention) so the values round-trip as text, exactly like
# the already text-safe values (e.g. "62.926.50", "0.0\u20830811").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.926.50'
$ws.Range("E2").Value = '  -1.46%  '

# Row 3
$ws.Range("D3").Value = '2.542.09'
$ws.Range("E3").Value = '  -0.32%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '''573.33'
$ws.Range("E5").Value = '  -0.52%  '

# Row 6
$ws.Range("D6").Value = '''145.71'
$ws.Range("E6").Value = '  -2.43%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.583'
$ws.Range("E8").Value = '  -1.51%  '

# Row 9
$ws.Range("D9").Value = '2.538.97'
$ws.Range("E9").Value = '  -0.42%  '

# Row 10
$ws.Range("E10").Value = '  -2.30%  '

# Row 11
$ws.Range("D11").Value = '''5.49'
$ws.Range("E11").Value = '  -5.18%  '

# Row 12
$ws.Range("E12").Value = '  -0.59%  '

# Row 13
$ws.Range("D13").Value = '''0.353'
$ws.Range("E13").Value = '  -1.82%  '

# Row 14
$ws.Range("D14").Value = '''27.27'
$ws.Range("E14").Value = '  -2.75%  '

# Row 15
$ws.Range("D15").Value = '2.995.07'
$ws.Range("E15").Value = '  -0.40%  '

# Row 16
$ws.Range("D16").Value = '62.833.25'
$ws.Range("E16").Value = '  -1.31%  '

# Row 17
$ws.Range("E17").Value = '  -1.92%  '

# Row 18
$ws.Range("D18").Value = '2.519.32'
$ws.Range("E18").Value = '  -1.18%  '

# Row 19
$ws.Range("D19").Value = '''11.28'
$ws.Range("E19").Value = '  -2.89%  '

# Row 20
$ws.Range("D20").Value = '''334.86'
$ws.Range("E20").Value = '  -2.97%  '

# Row 21
$ws.Range("D21").Value = '''4.32'
$ws.Range("E21").Value = '  -1.40%  '

# Row 22
$ws.Range("D22").Value = '''6.73'
$ws.Range("E22").Value = '  -2.64%  '

# Row 23
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").Value = '''65.14'
$ws.Range("E24").Value = '  -1.64%  '

# Row 25
$ws.Range("D25").Value = '''0.170'
$ws.Range("E25").Value = '  -1.02%  '

# Row 26
$ws.Range("D26").Value = '''1.59'
$ws.Range("E26").Value = '  +0.42%  '

# Row 27
$ws.Range("E27").Value = '  -0.14%  '

# Row 28
$ws.Range("D28").Value = '''8.32'
$ws.Range("E28").Value = '  -0.40%  '

# Row 29
$ws.Range("E29").Value = '  +1.13%  '

# Row 30
$ws.Range("D30").Value = '''7.23'
$ws.Range("E30").Value = '  +5.06%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0811'
$ws.Range("E31").Value = '  -3.33%  '

# Row 32
$ws.Range("D32").Value = '''1.86'
$ws.Range("E32").Value = '  -1.69%  '

# Row 33
$ws.Range("D33").Value = '''177.66'
$ws.Range("E33").Value = '  +0.52%  '

# Row 34
$ws.Range("D34").Value = '''1.53'
$ws.Range("E34").Value = '  -3.96%  '

# Row 35
$ws.Range("D35").Value = '''399.61'
$ws.Range("E35").Value = '  -5.92%  '

# Row 36
$ws.Range("D36").Value = '''19.04'
$ws.Range("E36").Value = '  -1.00%  '

# Row 37
$ws.Range("D37").Value = '''0.395'
$ws.Range("E37").Value = '  -2.76%  '

# Row 38
$ws.Range("E38").Value = '  +0.01%  '

# Row 39
$ws.Range("D39").Value = '''4.32'
$ws.Range("E39").Value = '  -3.23%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''1.73'
$ws.Range("E40").Value = '  -1.65%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("D42").Value = '''39.24'
$ws.Range("E42").Value = '  -3.31%  '

# Row 43
$ws.Range("D43").Value = '''150.42'
$ws.Range("E43").Value = '  -2.62%  '

# Row 44
$ws.Range("D44").Value = '''3.73'
$ws.Range("E44").Value = '  -2.09%  '

# Row 45
$ws.Range("D45").Value = '''20.68'
$ws.Range("E45").Value = '  -1.92%  '

# Row 46
$ws.Range("D46").Value = '''0.0532'
$ws.Range("E46").Value = '  -0.44%  '

# Row 47
$ws.Range("D47").Value = '''0.598'
$ws.Range("E47").Value = '  -2.68%  '

# Row 48
$ws.Range("D48").Value = '''0.0960'
$ws.Range("E48").Value = '  -0.97%  '

# Row 49
$ws.Range("D49").Value = '''0.0237'
$ws.Range("E49").Value = '  +1.82%  '

# Row 50
$ws.Range("D50").Value = '''18.11'
$ws.Range("E50").Value = '  -5.30%  '

# Row 51
$ws.Range("D51").Value = '''11.31'
$ws.Range("E51").Value = '  +0.42%  '
